# Automatic map update: remove the row for case -404 (Amenabar 3048),
# which shifts all subsequent rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(50).Delete()
